$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (URUGUAY - PRIMERA DIVISION, Fenix vs Racing Montevideo) was removed from
# the sheet. Deleting it shifts the following row (USA - MLS, New York Red Bulls vs
# Columbus Crew, formerly row 19) up to become the new row 18, and the sheet
# dimension shrinks from A1:BD19 to A1:BD18 automatically.
$ws.Rows.Item(18).Delete()

# Updated odds values (post row-shift, using final row numbers).

$ws.Range("G4").Value = 1.85
$ws.Range("I4").Value = 4.1
$ws.Range("X4").Value = 8.5
$ws.Range("AH4").Value = 10
$ws.Range("AJ4").Value = 13

$ws.Range("G5").Value = 2.9
$ws.Range("I5").Value = 2.88
$ws.Range("J5").Value = 4
$ws.Range("L5").Value = 3.75
$ws.Range("Z5").Value = 34
$ws.Range("AK5").Value = 29
$ws.Range("AO5").Value = 21
$ws.Range("AX5").Value = 19

$ws.Range("G6").Value = 2.2
$ws.Range("H6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6
$ws.Range("W6").Value = 6.5
$ws.Range("Z6").Value = 21
$ws.Range("AA6").Value = 21
$ws.Range("AC6").Value = 7.5
$ws.Range("AH6").Value = 9
$ws.Range("AO6").Value = 13
$ws.Range("AP6").Value = 26
$ws.Range("AV6").Value = 67
$ws.Range("AW6").Value = 5
$ws.Range("AY6").Value = 34

$ws.Range("G8").Value = 1.48
$ws.Range("H8").Value = 4.2
$ws.Range("K8").Value = 2.5
$ws.Range("L8").Value = 5.5
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 15
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.25
$ws.Range("W8").Value = 8.5
$ws.Range("AC8").Value = 15
$ws.Range("AJ8").Value = 19

$ws.Range("G10").Value = 1.67
$ws.Range("H10").Value = 3.4
$ws.Range("I10").Value = 6
$ws.Range("N10").Value = 7.5
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.62
$ws.Range("AA10").Value = 15
$ws.Range("AC10").Value = 7.5
$ws.Range("AK10").Value = 67
$ws.Range("AQ10").Value = 29
$ws.Range("AW10").Value = 7

$ws.Range("G11").Value = 3.2
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 2.4
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 3
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 1.83
$ws.Range("Z11").Value = 34
$ws.Range("AB11").Value = 41
$ws.Range("AH11").Value = 7.5
$ws.Range("AI11").Value = 11
$ws.Range("AP11").Value = 29
$ws.Range("AW11").Value = 4.33
$ws.Range("AX11").Value = 13

$ws.Range("G13").Value = 1.13
$ws.Range("H13").Value = 10
$ws.Range("K13").Value = 3
$ws.Range("N13").Value = 19
$ws.Range("Q13").Value = 1.5
$ws.Range("R13").Value = 2.5
$ws.Range("S13").Value = 1.25
$ws.Range("T13").Value = 3.75
$ws.Range("W13").Value = 7.5
$ws.Range("AC13").Value = 19
$ws.Range("AD13").Value = 21
$ws.Range("AQ13").Value = 10
$ws.Range("AT13").Value = 3.75

$ws.Range("H16").Value = 8.5
$ws.Range("K16").Value = 3.1
$ws.Range("N16").Value = 19
$ws.Range("V16").Value = 1.54
$ws.Range("W16").Value = 9
$ws.Range("AZ16").Value = 451

$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 2.8
$ws.Range("L18").Value = 3.25
$ws.Range("U18").Value = 1.57
$ws.Range("V18").Value = 2.25
$ws.Range("X18").Value = 13
$ws.Range("Z18").Value = 23
